$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tuesday 10.4.18"
$ws.Range("D4").Value = "Git work + arduino play"
$ws.Range("B4").Value = "1830- 2030"
$ws.Range("C4").Value = 2

$ws.Range("C5").Select()
